$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Dll3"
$ws.Cells.Item(2, 3).Value = "Notch2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 1.305319666666667
$ws.Cells.Item(2, 8).Value = 3.915959
$ws.Cells.Item(2, 9).Value = 0.4586709810613518
$ws.Cells.Item(2, 10).Value = 0.4586709810613519
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 1.910418
$ws.Cells.Item(2, 14).Value = 5.731254
$ws.Cells.Item(2, 15).Value = 0.01809124304049503
$ws.Cells.Item(2, 16).Value = 0.01809124304049503
$ws.Cells.Item(2, 17).Value = 2.493706186954
$ws.Cells.Item(2, 18).Value = 22.443355682586
$ws.Cells.Item(2, 19).Value = 0.008297928194003209
$ws.Cells.Item(2, 20).Value = 0.008297928194003209

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Dll3"
$ws.Cells.Item(3, 3).Value = "Notch2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 1.305319666666667
$ws.Cells.Item(3, 8).Value = 3.915959
$ws.Cells.Item(3, 9).Value = 0.4586709810613518
$ws.Cells.Item(3, 10).Value = 0.4586709810613519
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 31.995262
$ws.Cells.Item(3, 14).Value = 95.985786
$ws.Cells.Item(3, 15).Value = 0.302988173785169
$ws.Cells.Item(3, 16).Value = 0.302988173785169
$ws.Cells.Item(3, 17).Value = 41.76404472875267
$ws.Cells.Item(3, 18).Value = 375.876402558774
$ws.Cells.Item(3, 19).Value = 0.1389718829200308
$ws.Cells.Item(3, 20).Value = 0.1389718829200308

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Dll3"
$ws.Cells.Item(4, 3).Value = "Notch2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 1.305319666666667
$ws.Cells.Item(4, 8).Value = 3.915959
$ws.Cells.Item(4, 9).Value = 0.4586709810613518
$ws.Cells.Item(4, 10).Value = 0.4586709810613519
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 37.858701
$ws.Cells.Item(4, 14).Value = 113.576103
$ws.Cells.Item(4, 15).Value = 0.3585136661130873
$ws.Cells.Item(4, 16).Value = 0.3585136661130873
$ws.Cells.Item(4, 17).Value = 49.41770696975299
$ws.Cells.Item(4, 18).Value = 444.7593627277769
$ws.Cells.Item(4, 19).Value = 0.1644398149599917
$ws.Cells.Item(4, 20).Value = 0.1644398149599917

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Dll3"
$ws.Cells.Item(5, 3).Value = "Notch2"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 1.305319666666667
$ws.Cells.Item(5, 8).Value = 3.915959
$ws.Cells.Item(5, 9).Value = 0.4586709810613518
$ws.Cells.Item(5, 10).Value = 0.4586709810613519
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 33.83466466666667
$ws.Cells.Item(5, 14).Value = 101.503994
$ws.Cells.Item(5, 15).Value = 0.3204069170612486
$ws.Cells.Item(5, 16).Value = 0.3204069170612486
$ws.Cells.Item(5, 17).Value = 44.16505320447178
$ws.Cells.Item(5, 18).Value = 397.485478840246
$ws.Cells.Item(5, 19).Value = 0.1469613549873261
$ws.Cells.Item(5, 20).Value = 0.1469613549873261

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Dll3"
$ws.Cells.Item(6, 3).Value = "Notch2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.983774
$ws.Cells.Item(6, 8).Value = 2.951322
$ws.Cells.Item(6, 9).Value = 0.3456843539904149
$ws.Cells.Item(6, 10).Value = 0.3456843539904149
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.910418
$ws.Cells.Item(6, 14).Value = 5.731254
$ws.Cells.Item(6, 15).Value = 0.01809124304049503
$ws.Cells.Item(6, 16).Value = 0.01809124304049503
$ws.Cells.Item(6, 17).Value = 1.879419557532
$ws.Cells.Item(6, 18).Value = 16.914776017788
$ws.Cells.Item(6, 19).Value = 0.006253859663337112
$ws.Cells.Item(6, 20).Value = 0.006253859663337113

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Dll3"
$ws.Cells.Item(7, 3).Value = "Notch2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.983774
$ws.Cells.Item(7, 8).Value = 2.951322
$ws.Cells.Item(7, 9).Value = 0.3456843539904149
$ws.Cells.Item(7, 10).Value = 0.3456843539904149
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 31.995262
$ws.Cells.Item(7, 14).Value = 95.985786
$ws.Cells.Item(7, 15).Value = 0.302988173785169
$ws.Cells.Item(7, 16).Value = 0.302988173785169
$ws.Cells.Item(7, 17).Value = 31.476106878788
$ws.Cells.Item(7, 18).Value = 283.284961909092
$ws.Cells.Item(7, 19).Value = 0.1047382711216617
$ws.Cells.Item(7, 20).Value = 0.1047382711216617

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Dll3"
$ws.Cells.Item(8, 3).Value = "Notch2"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.983774
$ws.Cells.Item(8, 8).Value = 2.951322
$ws.Cells.Item(8, 9).Value = 0.3456843539904149
$ws.Cells.Item(8, 10).Value = 0.3456843539904149
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 37.858701
$ws.Cells.Item(8, 14).Value = 113.576103
$ws.Cells.Item(8, 15).Value = 0.3585136661130873
$ws.Cells.Item(8, 16).Value = 0.3585136661130873
$ws.Cells.Item(8, 17).Value = 37.244405717574
$ws.Cells.Item(8, 18).Value = 335.199651458166
$ws.Cells.Item(8, 19).Value = 0.1239325650670379
$ws.Cells.Item(8, 20).Value = 0.1239325650670379

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Dll3"
$ws.Cells.Item(9, 3).Value = "Notch2"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.983774
$ws.Cells.Item(9, 8).Value = 2.951322
$ws.Cells.Item(9, 9).Value = 0.3456843539904149
$ws.Cells.Item(9, 10).Value = 0.3456843539904149
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 33.83466466666667
$ws.Cells.Item(9, 14).Value = 101.503994
$ws.Cells.Item(9, 15).Value = 0.3204069170612486
$ws.Cells.Item(9, 16).Value = 0.3204069170612486
$ws.Cells.Item(9, 17).Value = 33.28566339778533
$ws.Cells.Item(9, 18).Value = 299.570970580068
$ws.Cells.Item(9, 19).Value = 0.1107596581383782
$ws.Cells.Item(9, 20).Value = 0.1107596581383782

$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Dll3"
$ws.Cells.Item(10, 3).Value = "Notch2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.1949326666666667
$ws.Cells.Item(10, 8).Value = 0.584798
$ws.Cells.Item(10, 9).Value = 0.0684965987597716
$ws.Cells.Item(10, 10).Value = 0.06849659875977161
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 1.910418
$ws.Cells.Item(10, 14).Value = 5.731254
$ws.Cells.Item(10, 15).Value = 0.01809124304049503
$ws.Cells.Item(10, 16).Value = 0.01809124304049503
$ws.Cells.Item(10, 17).Value = 0.372402875188
$ws.Cells.Item(10, 18).Value = 3.351625876692
$ws.Cells.Item(10, 19).Value = 0.001239188615610298
$ws.Cells.Item(10, 20).Value = 0.001239188615610298

$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Dll3"
$ws.Cells.Item(11, 3).Value = "Notch2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.1949326666666667
$ws.Cells.Item(11, 8).Value = 0.584798
$ws.Cells.Item(11, 9).Value = 0.0684965987597716
$ws.Cells.Item(11, 10).Value = 0.06849659875977161
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 31.995262
$ws.Cells.Item(11, 14).Value = 95.985786
$ws.Cells.Item(11, 15).Value = 0.302988173785169
$ws.Cells.Item(11, 16).Value = 0.302988173785169
$ws.Cells.Item(11, 17).Value = 6.236921742358667
$ws.Cells.Item(11, 18).Value = 56.13229568122801
$ws.Cells.Item(11, 19).Value = 0.02075365936871867
$ws.Cells.Item(11, 20).Value = 0.02075365936871867

$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Dll3"
$ws.Cells.Item(12, 3).Value = "Notch2"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.1949326666666667
$ws.Cells.Item(12, 8).Value = 0.584798
$ws.Cells.Item(12, 9).Value = 0.0684965987597716
$ws.Cells.Item(12, 10).Value = 0.06849659875977161
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 37.858701
$ws.Cells.Item(12, 14).Value = 113.576103
$ws.Cells.Item(12, 15).Value = 0.3585136661130873
$ws.Cells.Item(12, 16).Value = 0.3585136661130873
$ws.Cells.Item(12, 17).Value = 7.379897542466
$ws.Cells.Item(12, 18).Value = 66.419077882194
$ws.Cells.Item(12, 19).Value = 0.02455696673764287
$ws.Cells.Item(12, 20).Value = 0.02455696673764287

$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Dll3"
$ws.Cells.Item(13, 3).Value = "Notch2"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.1949326666666667
$ws.Cells.Item(13, 8).Value = 0.584798
$ws.Cells.Item(13, 9).Value = 0.0684965987597716
$ws.Cells.Item(13, 10).Value = 0.06849659875977161
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 33.83466466666667
$ws.Cells.Item(13, 14).Value = 101.503994
$ws.Cells.Item(13, 15).Value = 0.3204069170612486
$ws.Cells.Item(13, 16).Value = 0.3204069170612486
$ws.Cells.Item(13, 17).Value = 6.595481409245778
$ws.Cells.Item(13, 18).Value = 59.35933268321201
$ws.Cells.Item(13, 19).Value = 0.02194678403779976
$ws.Cells.Item(13, 20).Value = 0.02194678403779976

$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Dll3"
$ws.Cells.Item(14, 3).Value = "Notch2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.3618473333333334
$ws.Cells.Item(14, 8).Value = 1.085542
$ws.Cells.Item(14, 9).Value = 0.1271480661884616
$ws.Cells.Item(14, 10).Value = 0.1271480661884616
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 1.910418
$ws.Cells.Item(14, 14).Value = 5.731254
$ws.Cells.Item(14, 15).Value = 0.01809124304049503
$ws.Cells.Item(14, 16).Value = 0.01809124304049503
$ws.Cells.Item(14, 17).Value = 0.691279658852
$ws.Cells.Item(14, 18).Value = 6.221516929668
$ws.Cells.Item(14, 19).Value = 0.002300266567544408
$ws.Cells.Item(14, 20).Value = 0.002300266567544408

$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Dll3"
$ws.Cells.Item(15, 3).Value = "Notch2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.3618473333333334
$ws.Cells.Item(15, 8).Value = 1.085542
$ws.Cells.Item(15, 9).Value = 0.1271480661884616
$ws.Cells.Item(15, 10).Value = 0.1271480661884616
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 31.995262
$ws.Cells.Item(15, 14).Value = 95.985786
$ws.Cells.Item(15, 15).Value = 0.302988173785169
$ws.Cells.Item(15, 16).Value = 0.302988173785169
$ws.Cells.Item(15, 17).Value = 11.57740023400133
$ws.Cells.Item(15, 18).Value = 104.196602106012
$ws.Cells.Item(15, 19).Value = 0.03852436037475779
$ws.Cells.Item(15, 20).Value = 0.03852436037475779

$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Dll3"
$ws.Cells.Item(16, 3).Value = "Notch2"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.3618473333333334
$ws.Cells.Item(16, 8).Value = 1.085542
$ws.Cells.Item(16, 9).Value = 0.1271480661884616
$ws.Cells.Item(16, 10).Value = 0.1271480661884616
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 37.858701
$ws.Cells.Item(16, 14).Value = 113.576103
$ws.Cells.Item(16, 15).Value = 0.3585136661130873
$ws.Cells.Item(16, 16).Value = 0.3585136661130873
$ws.Cells.Item(16, 17).Value = 13.699070000314
$ws.Cells.Item(16, 18).Value = 123.291630002826
$ws.Cells.Item(16, 19).Value = 0.04558431934841486
$ws.Cells.Item(16, 20).Value = 0.04558431934841486

$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Dll3"
$ws.Cells.Item(17, 3).Value = "Notch2"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 0.3618473333333334
$ws.Cells.Item(17, 8).Value = 1.085542
$ws.Cells.Item(17, 9).Value = 0.1271480661884616
$ws.Cells.Item(17, 10).Value = 0.1271480661884616
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 33.83466466666667
$ws.Cells.Item(17, 14).Value = 101.503994
$ws.Cells.Item(17, 15).Value = 0.3204069170612486
$ws.Cells.Item(17, 16).Value = 0.3204069170612486
$ws.Cells.Item(17, 17).Value = 12.24298318386089
$ws.Cells.Item(17, 18).Value = 110.186848654748
$ws.Cells.Item(17, 19).Value = 0.04073911989774457
$ws.Cells.Item(17, 20).Value = 0.04073911989774457
